# "Doing Updates for Financials" - refresh the yearly financial figures
# on the SNVFF sheet with the latest reported numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SNVFF")

$ws.Range("D14").Value = 4200
$ws.Range("G14").Value = 10200
$ws.Range("D17").Value = 5600
$ws.Range("G17").Value = 11700
$ws.Range("H17").Value = 3100
$ws.Range("I17").Value = 3200
$ws.Range("J17").Value = 6800
$ws.Range("G18").Value = -11700
$ws.Range("H18").Value = -3100
$ws.Range("I18").Value = -3200
$ws.Range("J18").Value = -6800
$ws.Range("H21").Value = -3800
$ws.Range("D23").Value = -5600
$ws.Range("G23").Value = -12000
$ws.Range("H23").Value = -3900
$ws.Range("I23").Value = -3200
$ws.Range("J23").Value = -6800
$ws.Range("D26").Value = -5600
$ws.Range("G26").Value = -12000
$ws.Range("H26").Value = -3900
$ws.Range("I26").Value = -3200
$ws.Range("J26").Value = -6800
$ws.Range("D27").Value = -5600
$ws.Range("G27").Value = -12000
$ws.Range("H27").Value = -3900
$ws.Range("I27").Value = -3200
$ws.Range("J27").Value = -6700
$ws.Range("D33").Value = -5600
$ws.Range("G33").Value = -12000
$ws.Range("H33").Value = -3900
$ws.Range("I33").Value = -3200
$ws.Range("J33").Value = -6700
$ws.Range("D35").Value = -5600
$ws.Range("G35").Value = -12000
$ws.Range("H35").Value = -3900
$ws.Range("I35").Value = -3200
$ws.Range("J35").Value = -6700
$ws.Range("H41").Value = 1600
$ws.Range("I41").Value = 2700
$ws.Range("J43").Value = 500
$ws.Range("H46").Value = 1800
$ws.Range("I46").Value = 3700
$ws.Range("J46").Value = 2100
$ws.Range("H47").Value = 600
$ws.Range("I47").Value = 700
$ws.Range("H48").Value = 10800
$ws.Range("I48").Value = 10000
$ws.Range("J48").Value = 12300
$ws.Range("F54").Value = 0
$ws.Range("H54").Value = 13400
$ws.Range("I54").Value = 14500
$ws.Range("J54").Value = 14600
$ws.Range("E57").Value = 400
$ws.Range("I57").Value = 500
$ws.Range("D60").Value = 1200
$ws.Range("E60").Value = 400
$ws.Range("I60").Value = 500
$ws.Range("D66").Value = 1200
$ws.Range("E66").Value = 400
$ws.Range("I66").Value = 500
$ws.Range("D72").Value = -62200
$ws.Range("E72").Value = -56900
$ws.Range("F72").Value = -56000
$ws.Range("G72").Value = -55400
$ws.Range("H72").Value = -43400
$ws.Range("I72").Value = -39600
$ws.Range("J72").Value = -36400
$ws.Range("D76").Value = -1100
$ws.Range("E76").Value = 1500
$ws.Range("F76").Value = -300
$ws.Range("G76").Value = 200
$ws.Range("H76").Value = 12700
$ws.Range("I76").Value = 14000
$ws.Range("J76").Value = 13900
$ws.Range("D81").Value = -5600
$ws.Range("G81").Value = -12000
$ws.Range("H81").Value = -3900
$ws.Range("I81").Value = -3200
$ws.Range("J81").Value = -6700
$ws.Range("G89").Value = -1700
$ws.Range("H89").Value = -2800
$ws.Range("I89").Value = -2800
$ws.Range("J89").Value = -4800
$ws.Range("J91").Value = -4100
$ws.Range("G94").Value = 500
$ws.Range("J94").Value = -4200
$ws.Range("D100").Value = 2500
$ws.Range("I100").Value = 2900
$ws.Range("J100").Value = 7300
$ws.Range("G102").Value = -1200
$ws.Range("J102").Value = -1600
